# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / heading / "Outstanding" columns one
# place to the right, then make that sheet the active tab with L17 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N (existing N/O/P/.. shift right to O/P/Q/..)
$ws.Columns("N").Insert()

# Give the newly inserted column a width similar to its neighbours
$ws.Columns("N").ColumnWidth = 9.75

# Make "Repayment schedule" the active sheet/tab and select L17 on it
$ws.Activate()
$null = $ws.Range("L17").Select()
